# Scheduled runner: refresh market-board price snapshots (currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ) and the dependent Leve profit
# columns (LevePriceNQ/HQ, LeveProfitNQ/HQ) across all job sheets.
#
# Columns per sheet:
#   H = currentAveragePrice        I = currentAveragePriceNQ
#   J = currentAveragePriceHQ      K = LevePriceNQ
#   L = LevePriceHQ                M = LeveProfitNQ
#   N = LeveProfitHQ
#
# All cells here are static snapshot values (no formulas in this workbook),
# so each updated figure is written as a literal value. Where a profit
# column's new value is immaterial (no HQ/NQ leve price applies) the cell is
# cleared entirely rather than left with a stale number.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 2099.5  # H28: was 2199.7144
$ws.Cells.Item(28, 9).Value = 1970.8572  # I28: was 2066.3333
$ws.Cells.Item(28, 11).Value = 1970.8572  # K28: was 2066.3333
$ws.Cells.Item(28, 13).Value = -1485.8572  # M28: was -1581.3333
$ws.Cells.Item(74, 8).Value = 2497.5  # H74: was 0
$ws.Cells.Item(74, 9).Value = 2497.5  # I74: was 0
$ws.Cells.Item(74, 11).Value = 2497.5  # K74: was 0
$ws.Cells.Item(74, 13).Value = -1561.5  # M74: was None
$ws.Cells.Item(77, 8).Value = 2497.5  # H77: was 0
$ws.Cells.Item(77, 9).Value = 2497.5  # I77: was 0
$ws.Cells.Item(77, 11).Value = 12487.5  # K77: was 0
$ws.Cells.Item(77, 13).Value = -7807.5  # M77: was None
$ws.Cells.Item(86, 8).Value = 2818.8  # H86: was 3486
$ws.Cells.Item(86, 9).Value = 2818.8  # I86: was 3486
$ws.Cells.Item(86, 11).Value = 2818.8  # K86: was 3486
$ws.Cells.Item(86, 13).Value = -1695.8  # M86: was -2363
$ws.Cells.Item(89, 8).Value = 2818.8  # H89: was 3486
$ws.Cells.Item(89, 9).Value = 2818.8  # I89: was 3486
$ws.Cells.Item(89, 11).Value = 14094  # K89: was 17430
$ws.Cells.Item(89, 13).Value = -8478  # M89: was -11814
$ws.Cells.Item(129, 8).Value = 5000  # H129: was 3450
$ws.Cells.Item(129, 9).Value = 0  # I129: was 350
$ws.Cells.Item(129, 11).Value = 0  # K129: was 1050
$ws.Cells.Item(129, 13).ClearContents()  # M129: was 3950

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(39, 8).Value = 12500  # H39: was 15000
$ws.Cells.Item(39, 9).Value = 12500  # I39: was 15000
$ws.Cells.Item(39, 11).Value = 12500  # K39: was 15000
$ws.Cells.Item(39, 13).Value = -11980  # M39: was -14480
$ws.Cells.Item(50, 8).Value = 14002.1  # H50: was 13330.333
$ws.Cells.Item(50, 9).Value = 6672.6665  # I50: was 3997.6
$ws.Cells.Item(50, 11).Value = 6672.6665  # K50: was 3997.6
$ws.Cells.Item(50, 13).Value = -5958.6665  # M50: was -3283.6
$ws.Cells.Item(122, 8).Value = 1300  # H122: was 1637.5
$ws.Cells.Item(122, 9).Value = 1210  # I122: was 1600
$ws.Cells.Item(122, 11).Value = 3630  # K122: was 4800
$ws.Cells.Item(122, 13).Value = -1180  # M122: was -2350

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(11, 8).Value = 1692.875  # H11: was 3005.5
$ws.Cells.Item(11, 9).Value = 137.75  # I11: was 20
$ws.Cells.Item(11, 10).Value = 3248  # J11: was 4000.6667
$ws.Cells.Item(11, 11).Value = 137.75  # K11: was 20
$ws.Cells.Item(11, 12).Value = 3248  # L11: was 4000.6667
$ws.Cells.Item(11, 13).Value = 2.25  # M11: was 120
$ws.Cells.Item(11, 14).Value = -3528  # N11: was -4280.6667
$ws.Cells.Item(19, 8).Value = 9594.444  # H19: was 9319.9
$ws.Cells.Item(19, 10).Value = 12770  # J19: was 11783.167
$ws.Cells.Item(19, 12).Value = 12770  # L19: was 11783.167
$ws.Cells.Item(19, 14).Value = -13116  # N19: was -12129.167
$ws.Cells.Item(99, 8).Value = 5224.5  # H99: was 6500.25
$ws.Cells.Item(99, 9).Value = 5960  # I99: was 7333.3335
$ws.Cells.Item(99, 10).Value = 3998.6667  # J99: was 4001
$ws.Cells.Item(99, 11).Value = 5960  # K99: was 7333.3335
$ws.Cells.Item(99, 12).Value = 3998.6667  # L99: was 4001
$ws.Cells.Item(99, 13).Value = -4462  # M99: was -5835.3335
$ws.Cells.Item(99, 14).Value = -6994.6667  # N99: was -6997

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(2, 8).Value = 1033.8889  # H2: was 1053.3334
$ws.Cells.Item(2, 9).Value = 845  # I2: was 781.4286
$ws.Cells.Item(2, 10).Value = 1411.6666  # J2: was 2005
$ws.Cells.Item(2, 11).Value = 845  # K2: was 781.4286
$ws.Cells.Item(2, 12).Value = 1411.6666  # L2: was 2005
$ws.Cells.Item(2, 13).Value = -732  # M2: was -668.4286
$ws.Cells.Item(2, 14).Value = -1637.6666  # N2: was -2231
$ws.Cells.Item(5, 8).Value = 2376.3333  # H5: was 4241.25
$ws.Cells.Item(5, 9).Value = 2654.8333  # I5: was 5978.5
$ws.Cells.Item(5, 10).Value = 1819.3334  # J5: was 2504
$ws.Cells.Item(5, 11).Value = 2654.8333  # K5: was 5978.5
$ws.Cells.Item(5, 12).Value = 1819.3334  # L5: was 2504
$ws.Cells.Item(5, 13).Value = -2542.8333  # M5: was -5866.5
$ws.Cells.Item(5, 14).Value = -2043.3334  # N5: was -2728
$ws.Cells.Item(7, 8).Value = 85.5  # H7: was 114.57143
$ws.Cells.Item(7, 9).Value = 44.285713  # I7: was 71.4
$ws.Cells.Item(7, 10).Value = 181.66667  # J7: was 222.5
$ws.Cells.Item(7, 11).Value = 44.285713  # K7: was 71.4
$ws.Cells.Item(7, 12).Value = 181.66667  # L7: was 222.5
$ws.Cells.Item(7, 13).Value = 68.714287  # M7: was 41.59999999999999
$ws.Cells.Item(7, 14).Value = -407.66667  # N7: was -448.5
$ws.Cells.Item(8, 8).Value = 8870  # H8: was 2499
$ws.Cells.Item(8, 9).Value = 0  # I8: was 9
$ws.Cells.Item(8, 10).Value = 8870  # J8: was 4989
$ws.Cells.Item(8, 11).Value = 0  # K8: was 9
$ws.Cells.Item(8, 12).Value = 8870  # L8: was 4989
$ws.Cells.Item(8, 13).ClearContents()  # M8: was 131
$ws.Cells.Item(8, 14).Value = -9150  # N8: was -5269
$ws.Cells.Item(10, 8).Value = 889.8  # H10: was 1858.4615
$ws.Cells.Item(10, 9).Value = 498.625  # I10: was 706.875
$ws.Cells.Item(10, 10).Value = 1336.8572  # J10: was 3701
$ws.Cells.Item(10, 11).Value = 498.625  # K10: was 706.875
$ws.Cells.Item(10, 12).Value = 1336.8572  # L10: was 3701
$ws.Cells.Item(10, 13).Value = -359.625  # M10: was -567.875
$ws.Cells.Item(10, 14).Value = -1614.8572  # N10: was -3979
$ws.Cells.Item(11, 8).Value = 10360  # H11: was 20753
$ws.Cells.Item(11, 9).Value = 10  # I11: was 2000
$ws.Cells.Item(11, 10).Value = 13810  # J11: was 27004
$ws.Cells.Item(11, 11).Value = 10  # K11: was 2000
$ws.Cells.Item(11, 12).Value = 13810  # L11: was 27004
$ws.Cells.Item(11, 13).Value = 130  # M11: was -1860
$ws.Cells.Item(11, 14).Value = -14090  # N11: was -27284
$ws.Cells.Item(12, 8).Value = 984.5455  # H12: was 1132.2727
$ws.Cells.Item(12, 9).Value = 450  # I12: was 675
$ws.Cells.Item(12, 10).Value = 1103.3334  # J12: was 1233.8889
$ws.Cells.Item(12, 11).Value = 450  # K12: was 675
$ws.Cells.Item(12, 12).Value = 1103.3334  # L12: was 1233.8889
$ws.Cells.Item(12, 13).Value = -280  # M12: was -505
$ws.Cells.Item(12, 14).Value = -1443.3334  # N12: was -1573.8889
$ws.Cells.Item(13, 8).Value = 9819.4  # H13: was 750
$ws.Cells.Item(13, 9).Value = 600  # I13: was 500
$ws.Cells.Item(13, 10).Value = 12124.25  # J13: was 1000
$ws.Cells.Item(13, 11).Value = 600  # K13: was 500
$ws.Cells.Item(13, 12).Value = 12124.25  # L13: was 1000
$ws.Cells.Item(13, 13).Value = -461  # M13: was -361
$ws.Cells.Item(13, 14).Value = -12402.25  # N13: was -1278
$ws.Cells.Item(14, 8).Value = 860  # H14: was 505
$ws.Cells.Item(14, 9).Value = 860  # I14: was 505
$ws.Cells.Item(14, 11).Value = 860  # K14: was 505
$ws.Cells.Item(14, 13).Value = -690  # M14: was -335
$ws.Cells.Item(19, 8).Value = 202.76923  # H19: was 174.33333
$ws.Cells.Item(19, 9).Value = 167.72728  # I19: was 174.33333
$ws.Cells.Item(19, 10).Value = 395.5  # J19: was 0
$ws.Cells.Item(19, 11).Value = 167.72728  # K19: was 174.33333
$ws.Cells.Item(19, 12).Value = 395.5  # L19: was 0
$ws.Cells.Item(19, 13).Value = 2.272719999999993  # M19: was -4.333329999999989
$ws.Cells.Item(19, 14).Value = -735.5  # N19: was None
$ws.Cells.Item(24, 8).Value = 202.76923  # H24: was 174.33333
$ws.Cells.Item(24, 9).Value = 167.72728  # I24: was 174.33333
$ws.Cells.Item(24, 10).Value = 395.5  # J24: was 0
$ws.Cells.Item(24, 11).Value = 167.72728  # K24: was 174.33333
$ws.Cells.Item(24, 12).Value = 395.5  # L24: was 0
$ws.Cells.Item(24, 13).Value = 2.272719999999993  # M24: was -4.333329999999989
$ws.Cells.Item(24, 14).Value = -735.5  # N24: was None
$ws.Cells.Item(35, 8).Value = 5644.75  # H35: was 5189.778
$ws.Cells.Item(35, 10).Value = 6515  # J35: was 5273.75
$ws.Cells.Item(35, 12).Value = 6515  # L35: was 5273.75
$ws.Cells.Item(35, 14).Value = -7103  # N35: was -5861.75

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 91.666664  # H11: was 100
$ws.Cells.Item(11, 9).Value = 91.666664  # I11: was 100
$ws.Cells.Item(11, 11).Value = 274.999992  # K11: was 300
$ws.Cells.Item(11, 13).Value = -134.999992  # M11: was -160
$ws.Cells.Item(75, 8).Value = 2650  # H75: was 983.1429
$ws.Cells.Item(75, 9).Value = 0  # I75: was 142.33333
$ws.Cells.Item(75, 10).Value = 2650  # J75: was 1613.75
$ws.Cells.Item(75, 11).Value = 0  # K75: was 426.99999
$ws.Cells.Item(75, 12).Value = 7950  # L75: was 4841.25
$ws.Cells.Item(75, 13).ClearContents()  # M75: was 571.00001
$ws.Cells.Item(75, 14).Value = -9946  # N75: was -6837.25
$ws.Cells.Item(78, 8).Value = 2650  # H78: was 983.1429
$ws.Cells.Item(78, 9).Value = 0  # I78: was 142.33333
$ws.Cells.Item(78, 10).Value = 2650  # J78: was 1613.75
$ws.Cells.Item(78, 11).Value = 0  # K78: was 1280.99997
$ws.Cells.Item(78, 12).Value = 23850  # L78: was 14523.75
$ws.Cells.Item(78, 13).ClearContents()  # M78: was 3711.00003
$ws.Cells.Item(78, 14).Value = -33834  # N78: was -24507.75
$ws.Cells.Item(86, 8).Value = 1581.4445  # H86: was 1955.6666
$ws.Cells.Item(86, 9).Value = 1870.4286  # I86: was 1936.8
$ws.Cells.Item(86, 10).Value = 570  # J86: was 2050
$ws.Cells.Item(86, 11).Value = 5611.2858  # K86: was 5810.4
$ws.Cells.Item(86, 12).Value = 1710  # L86: was 6150
$ws.Cells.Item(86, 13).Value = -4425.2858  # M86: was -4624.4
$ws.Cells.Item(86, 14).Value = -4082  # N86: was -8522
$ws.Cells.Item(89, 8).Value = 1581.4445  # H89: was 1955.6666
$ws.Cells.Item(89, 9).Value = 1870.4286  # I89: was 1936.8
$ws.Cells.Item(89, 10).Value = 570  # J89: was 2050
$ws.Cells.Item(89, 11).Value = 16833.8574  # K89: was 17431.2
$ws.Cells.Item(89, 12).Value = 5130  # L89: was 18450
$ws.Cells.Item(89, 13).Value = -10905.8574  # M89: was -11503.2
$ws.Cells.Item(89, 14).Value = -16986  # N89: was -30306
$ws.Cells.Item(92, 8).Value = 985  # H92: was 990
$ws.Cells.Item(92, 10).Value = 975  # J92: was 0
$ws.Cells.Item(92, 12).Value = 2925  # L92: was 0
$ws.Cells.Item(92, 14).Value = -5421  # N92: was None
$ws.Cells.Item(98, 8).Value = 1788.4546  # H98: was 1788.8182
$ws.Cells.Item(98, 9).Value = 1776.7142  # I98: was 1777.2858
$ws.Cells.Item(98, 11).Value = 5330.142599999999  # K98: was 5331.857400000001
$ws.Cells.Item(98, 13).Value = -3832.142599999999  # M98: was -3833.857400000001

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 0  # H80: was 4024.25
$ws.Cells.Item(80, 9).Value = 0  # I80: was 3399
$ws.Cells.Item(80, 10).Value = 0  # J80: was 4649.5
$ws.Cells.Item(80, 11).Value = 0  # K80: was 3399
$ws.Cells.Item(80, 12).Value = 0  # L80: was 4649.5
$ws.Cells.Item(80, 13).ClearContents()  # M80: was -2401
$ws.Cells.Item(80, 14).ClearContents()  # N80: was -6645.5
$ws.Cells.Item(83, 8).Value = 0  # H83: was 4024.25
$ws.Cells.Item(83, 9).Value = 0  # I83: was 3399
$ws.Cells.Item(83, 10).Value = 0  # J83: was 4649.5
$ws.Cells.Item(83, 11).Value = 0  # K83: was 16995
$ws.Cells.Item(83, 12).Value = 0  # L83: was 23247.5
$ws.Cells.Item(83, 13).ClearContents()  # M83: was -12003
$ws.Cells.Item(83, 14).ClearContents()  # N83: was -33231.5
$ws.Cells.Item(107, 8).Value = 1002  # H107: was 0
$ws.Cells.Item(107, 9).Value = 1002  # I107: was 0
$ws.Cells.Item(107, 11).Value = 1002  # K107: was 0
$ws.Cells.Item(107, 13).Value = 918  # M107: was None
$ws.Cells.Item(122, 8).Value = 0  # H122: was 4081.125
$ws.Cells.Item(122, 9).Value = 0  # I122: was 4437.5
$ws.Cells.Item(122, 10).Value = 0  # J122: was 3724.75
$ws.Cells.Item(122, 11).Value = 0  # K122: was 13312.5
$ws.Cells.Item(122, 12).Value = 0  # L122: was 11174.25
$ws.Cells.Item(122, 13).ClearContents()  # M122: was -10862.5
$ws.Cells.Item(122, 14).ClearContents()  # N122: was -16074.25

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 891.5  # H81: was 855.3333
$ws.Cells.Item(81, 10).Value = 1000  # J81: was 0
$ws.Cells.Item(81, 12).Value = 2000  # L81: was 0
$ws.Cells.Item(81, 14).Value = -4122  # N81: was None
$ws.Cells.Item(84, 8).Value = 891.5  # H84: was 855.3333
$ws.Cells.Item(84, 10).Value = 1000  # J84: was 0
$ws.Cells.Item(84, 12).Value = 10000  # L84: was 0
$ws.Cells.Item(84, 14).Value = -20608  # N84: was None